# Adds a new '2022-Q1' sheet (by renaming + repopulating the former
# '总计' sheet, which keeps its sheetId/rId) and appends a brand-new
# '总计' sheet at the end with the summary row for 2022-Q1 prepended.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count
$quarterSheet = $wb.Worksheets.Item($sheetCount)

# Helper: apply the bold / centered / thin-bordered "header style" used
# throughout this workbook for label columns + the index column.
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Helper: write a value that looks numeric ("68.43") but must be stored
# as literal text, matching the source data export, without leaving a
# quote-prefix style behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---- Step 1: the old "总计" sheet becomes "2022-Q1" --------------------
$quarterSheet.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $quarterSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    if ($col -ge 5) {
        Set-HeaderStyle $cell
    }
}

$quarterData = @(
    @("510810", "汇添富中证上海国企ETF", "68.43", "99.71", "3.79", "2.5935", 7),
    @("161721", "招商沪深300地产等权重指数", "9.97", "94.51", "13.39", "1.3350", 1),
    @("100032", "富国中证红利指数增强", "56.72", "94.48", "2.07", "1.1741", 6),
    @("512200", "南方中证全指房地产ETF", "28.63", "99.85", "2.64", "0.7558", 7),
    @("160218", "国泰国证房地产行业指数", "6.35", "95.04", "2.77", "0.1759", 6),
    @("530011", "建信内生动力混合", "3.11", "80.80", "4.30", "0.1337", 9),
    @("160628", "鹏华中证800地产指数（LOF）", "3.38", "94.35", "3.76", "0.1271", 7),
    @("515450", "南方标普中国A股大盘红利低波50ETF", "2.04", "99.51", "3.61", "0.0736", 5),
    @("001276", "建信新经济灵活配置混合", "1.59", "84.07", "4.22", "0.0671", 9),
    @("515060", "华夏中证全指房地产ETF", "2.37", "98.82", "2.62", "0.0621", 7),
    @("008114", "天弘中证红利低波动100指数A", "3.16", "92.60", "1.84", "0.0581", 4),
    @("009347", "中融价值成长6个月持有期混合A", "1.69", "93.13", "3.01", "0.0509", 10),
    @("008115", "天弘中证红利低波动100指数C", "2.37", "92.60", "1.84", "0.0436", 4),
    @("515100", "景顺长城中证红利低波动100ETF", "1.25", "97.96", "1.93", "0.0241", 5),
    @("512530", "建信沪深300红利ETF", "0.52", "97.34", "2.70", "0.0140", 6),
    @("002495", "前海开源量化优选灵活配置混合A", "0.49", "93.24", "2.81", "0.0138", 6),
    @("002496", "前海开源量化优选灵活配置混合C", "0.26", "93.24", "2.81", "0.0073", 6),
    @("510190", "华安上证龙头ETF", "0.53", "97.53", "1.32", "0.0070", 4),
    @("009348", "中融价值成长6个月持有期混合C", "0.21", "93.13", "3.01", "0.0063", 10),
    @("010404", "博道盛利6个月持有期混合", "1.29", "34.13", "0.42", "0.0054", 10),
    @("002334", "汇丰晋信大盘波动精选股票A", "0.17", "88.41", "1.99", "0.0034", 4),
    @("002335", "汇丰晋信大盘波动精选股票C", "0.02", "88.41", "1.99", "0.0004", 4)
)

for ($i = 0; $i -lt $quarterData.Count; $i++) {
    $r = $i + 2
    $row = $quarterData[$i]
    $aCell = $quarterSheet.Cells.Item($r, 1)
    $aCell.Value = $i
    if ($r -ge 7) {
        Set-HeaderStyle $aCell
    }
    $quarterSheet.Cells.Item($r, 2).Value = $row[0]
    $quarterSheet.Cells.Item($r, 3).Value = $row[1]
    Set-TextValue $quarterSheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $quarterSheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $quarterSheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $quarterSheet.Cells.Item($r, 7) $row[5]
    $quarterSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---- Step 2: brand-new "总计" sheet appended after "2022-Q1" -----------
$totalSheet = $wb.Worksheets.Add($null, $quarterSheet)
$totalSheet.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $cell = $totalSheet.Cells.Item(1, $col)
    $cell.Value = $totalHeaders[$col - 2]
    Set-HeaderStyle $cell
}

$totalData = @(
    @("2022-Q1", 22, 6.73),
    @("2021-Q4", 7, 4.99),
    @("2021-Q3", 10, 4.48),
    @("2021-Q2", 6, 5.07),
    @("2021-Q1", 9, 5.57),
    @("2020-Q4", 10, 6.03)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $r = $i + 2
    $row = $totalData[$i]
    $aCell = $totalSheet.Cells.Item($r, 1)
    $aCell.Value = $i
    Set-HeaderStyle $aCell
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
}

Write-Output "2022-Q1 sheet added; 总计 sheet rebuilt"
